# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
#
# Rows 9, 10, 11 get cyclically rotated:
#   new row 9  <- old row 11 (B:AD)
#   new row 10 <- old row 9  (B:AD)
#   new row 11 <- old row 10 (B:AD)
# Rows 131, 132 are swapped:
#   new row 131 <- old row 132 (B:AD)
#   new row 132 <- old row 131 (B:AD)
# Column A (the sequential row index) stays untouched on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the "before" values of the data columns (B:AD) for the rows involved ---
$row9  = $ws.Range("B9:AD9").Value2
$row10 = $ws.Range("B10:AD10").Value2
$row11 = $ws.Range("B11:AD11").Value2

$row131 = $ws.Range("B131:AD131").Value2
$row132 = $ws.Range("B132:AD132").Value2

# --- write back the rotated/swapped values ---
$ws.Range("B9:AD9").Value2   = $row11
$ws.Range("B10:AD10").Value2 = $row9
$ws.Range("B11:AD11").Value2 = $row10

$ws.Range("B131:AD131").Value2 = $row132
$ws.Range("B132:AD132").Value2 = $row131
